$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'" + "63.272.15"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.83%  "
$c = $ws.Range("D3")
$c.Value = "'" + "2.646.29"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.Value = "'" + "594.13"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  +0.05%  "
$c = $ws.Range("D8")
$c.Value = "'" + "0.585"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.90%  "
$c = $ws.Range("D9")
$c.Value = "'" + "2.647.66"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  +0.86%  "
$c = $ws.Range("D13")
$c.Value = "'" + "0.355"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.51%  "
$c = $ws.Range("D14")
$c.Value = "'" + "27.36"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.40%  "
$c = $ws.Range("D15")
$c.Value = "'" + "3.120.78"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.49%  "
$c = $ws.Range("D16")
$c.Value = "'" + "63.239.47"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  -0.59%  "
$c = $ws.Range("D18")
$c.Value = "'" + "2.625.23"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.41%  "
$c = $ws.Range("D19")
$c.Value = "'" + "11.39"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.37%  "
$c = $ws.Range("D20")
$c.Value = "'" + "338.95"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  -0.15%  "
$c = $ws.Range("D22")
$c.Value = "'" + "6.72"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("E23").Value = "  +0.09%  "
$c = $ws.Range("D24")
$c.Value = "'" + "66.95"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.78%  "
$c = $ws.Range("D25")
$c.Value = "'" + "1.66"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.78%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D26")
$c.Value = "'" + "1.52"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D27")
$c.Value = "'" + "0.165"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  +0.23%  "
$c = $ws.Range("D29")
$c.Value = "'" + "8.39"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("E30").Value = "  -2.05%  "
$c = $ws.Range("D31")
$c.Value = "'" + "525.88"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +15.81%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D32")
$c.Value = "'" + "1.81"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +11.26%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D33")
$c.Value = "'" + "1.96"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  -1.08%  "
$c = $ws.Range("D35")
$c.Value = "'" + "174.47"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "
$c = $ws.Range("D36")
$c.Value = "'" + "4.88"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +9.73%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +0.47%  "
$c = $ws.Range("D39")
$c.Value = "'" + "19.00"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  +6.27%  "
$c = $ws.Range("D41")
$c.Value = "'" + "172.06"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.16%  "
$c = $ws.Range("D42")
$c.Value = "'" + "0.999"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "
$c = $ws.Range("D43")
$c.Value = "'" + "40.22"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.98%  "
$c = $ws.Range("D44")
$c.Value = "'" + "3.72"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "
$c = $ws.Range("D45")
$c.Value = "'" + "21.99"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.17%  "
$c = $ws.Range("D46")
$c.Value = "'" + "0.0559"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.30%  "
$c = $ws.Range("D47")
$c.Value = "'" + "0.631"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "
$c = $ws.Range("D48")
$c.Value = "'" + "0.0959"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  +1.32%  "
$c = $ws.Range("D50")
$c.Value = "'" + "18.53"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  -0.65%  "
